# Implemented check for repeated activity names.
# The "Workflow" checklist sheet previously had all checks enabled ("Yes")
# except for "Variable scope is not the innermost", "Repeated display names
# for activities" and "Unused variables" which were disabled ("No").
# This change flips most checks to disabled ("No") and enables the
# "Repeated display names for activities" check ("Yes") with a threshold
# argument of 1 in column D.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Workflow")

# Disable checks in rows 2-18 and 20 (column A: Yes -> No)
$ws.Range("A2").Value = "No"
$ws.Range("A3").Value = "No"
$ws.Range("A4").Value = "No"
$ws.Range("A5").Value = "No"
$ws.Range("A6").Value = "No"
$ws.Range("A7").Value = "No"
$ws.Range("A8").Value = "No"
$ws.Range("A9").Value = "No"
$ws.Range("A10").Value = "No"
$ws.Range("A11").Value = "No"
$ws.Range("A12").Value = "No"
$ws.Range("A14").Value = "No"
$ws.Range("A15").Value = "No"
$ws.Range("A16").Value = "No"
$ws.Range("A17").Value = "No"
$ws.Range("A18").Value = "No"
$ws.Range("A20").Value = "No"

# Enable the "Repeated display names for activities" check (row 19) and
# set its threshold argument to 1.
$ws.Range("A19").Value = "Yes"
$ws.Range("D19").Value = 1

# Update the saved selection/scroll position to row 19.
$ws.Activate()
$win = $excel.ActiveWindow
$win.ScrollRow = 19
$win.ScrollColumn = 1
$ws.Range("E19").Select()
